$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174255728721619
$ws.Range("B1").Value = 2.394660234451294
$ws.Range("D1").Value = 2.359913110733032
$ws.Range("E1").Value = 1.20861828327179
